{"js": "// Extended credit history calculation on issues more than 1.5M\n// 1) Drop the now-unused DDE bookmark around the \"additionally contract\n//    guarantee issue with cost\" placeholder paragraph.\n// 2) Add two new DDE bookmarks around the \"\u041a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u0441\u043b\u0443\u0447\u0430\u0435\u0432 \u043f\u0440\u043e\u0441\u0440\u043e\u0447\u043a\u0438\"\n//    and \"\u0421\u043e\u0432\u043e\u043a\u0443\u043f\u043d\u043e\u0435 \u043a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u0434\u043d\u0435\u0439 \u043f\u0440\u043e\u0441\u0440\u043e\u0447\u043a\u0438\" labels in the credit\n//    history table.\n// 3) Replace the placeholder hyphen-minus \"-\" with an em dash \"\u2014\" in the\n//    two corresponding value cells.\n\nconst body = context.document.body;\n\n// --- 1) remove obsolete bookmark -------------------------------------------\ncontext.document.deleteBookmark(\"__DdeLink__1171_85680802\");\n\n// --- 2) add bookmarks around the two credit-history labels ------------------\nconst label1Results = body.search(\"\u041a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u0441\u043b\u0443\u0447\u0430\u0435\u0432 \u043f\u0440\u043e\u0441\u0440\u043e\u0447\u043a\u0438\", { matchCase: true });\nconst label2Results = body.search(\"\u0421\u043e\u0432\u043e\u043a\u0443\u043f\u043d\u043e\u0435 \u043a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u0434\u043d\u0435\u0439 \u043f\u0440\u043e\u0441\u0440\u043e\u0447\u043a\u0438\", { matchCase: true });\nlabel1Results.load(\"items\");\nlabel2Results.load(\"items\");\nawait context.sync();\n\nlabel1Results.items[0].insertBookmark(\"__DdeLink__6205_85680802\");\nlabel2Results.items[0].insertBookmark(\"__DdeLink__6207_85680802\");\nawait context.sync();\n\n// --- 3) replace the placeholder \"-\" with an em dash in both value cells -----\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => {\n  if (p.text === \"-\") {\n    p.getRange().insertText(\"\u2014\", Word.InsertLocation.replace);\n  }\n});\nawait context.sync();\n", "ps1": "# Extended credit history calculation on issues more than 1.5M\n# 1) Drop the now-unused DDE bookmark around the \"additionally contract\n#    guarantee issue with cost\" placeholder paragraph.\n# 2) Add two new DDE bookmarks around the \"\u041a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u0441\u043b\u0443\u0447\u0430\u0435\u0432 \u043f\u0440\u043e\u0441\u0440\u043e\u0447\u043a\u0438\"\n#    and \"\u0421\u043e\u0432\u043e\u043a\u0443\u043f\u043d\u043e\u0435 \u043a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u0434\u043d\u0435\u0439 \u043f\u0440\u043e\u0441\u0440\u043e\u0447\u043a\u0438\" labels in the credit\n#    history table.\n# 3) Replace the placeholder hyphen-minus \"-\" with an em dash \"\u2014\" in the\n#    two corresponding value cells.\n\n$d = $word.ActiveDocument\n\n# --- 1) remove obsolete bookmark -------------------------------------------------\n$oldBookmarkName = \"__DdeLink__1171_85680802\"\nif ($d.Bookmarks.Exists($oldBookmarkName)) {\n    $d.Bookmarks($oldBookmarkName).Delete()\n}\n\n# --- 2) add bookmark around \"\u041a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u0441\u043b\u0443\u0447\u0430\u0435\u0432 \u043f\u0440\u043e\u0441\u0440\u043e\u0447\u043a\u0438\" ----------------------\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Execute(\"\u041a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u0441\u043b\u0443\u0447\u0430\u0435\u0432 \u043f\u0440\u043e\u0441\u0440\u043e\u0447\u043a\u0438\") | Out-Null\n$d.Bookmarks.Add(\"__DdeLink__6205_85680802\", $rng1)\n\n# --- 2) add bookmark around \"\u0421\u043e\u0432\u043e\u043a\u0443\u043f\u043d\u043e\u0435 \u043a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u0434\u043d\u0435\u0439 \u043f\u0440\u043e\u0441\u0440\u043e\u0447\u043a\u0438\" --------------\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Execute(\"\u0421\u043e\u0432\u043e\u043a\u0443\u043f\u043d\u043e\u0435 \u043a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u0434\u043d\u0435\u0439 \u043f\u0440\u043e\u0441\u0440\u043e\u0447\u043a\u0438\") | Out-Null\n$d.Bookmarks.Add(\"__DdeLink__6207_85680802\", $rng2)\n\n# --- 3) replace the placeholder \"-\" with an em dash in both value cells ---------\n# (the value cell immediately follows its label cell in document order, so we\n#  search for the next standalone \"-\" right after each label)\nforeach ($labelText in @(\"\u041a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u0441\u043b\u0443\u0447\u0430\u0435\u0432 \u043f\u0440\u043e\u0441\u0440\u043e\u0447\u043a\u0438\", \"\u0421\u043e\u0432\u043e\u043a\u0443\u043f\u043d\u043e\u0435 \u043a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u0434\u043d\u0435\u0439 \u043f\u0440\u043e\u0441\u0440\u043e\u0447\u043a\u0438\")) {\n    $lr = $d.Content\n    $lr.Find.ClearFormatting()\n    $lr.Find.Execute($labelText) | Out-Null\n    $valueRange = $d.Range($lr.End, $lr.End + 40)\n    $valueRange.Find.ClearFormatting()\n    $valueRange.Find.Execute(\"-\") | Out-Null\n    if ($valueRange.Text -eq \"-\") {\n        $valueRange.Text = \"\u2014\"\n    }\n}\n"}
